$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.229.02"
$ws.Range("E2").Value = "  +0.81%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.854.49"
$ws.Range("E3").Value = "  +1.35%  "
$ws.Range("E4").Value = "  -0.34%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.76"
$ws.Range("E5").Value = "  +0.74%  "
$ws.Range("E6").Value = "  -0.27%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4644"
$ws.Range("E7").Value = "  +0.37%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07291"
$ws.Range("E9").Value = "  -0.71%  "
$ws.Range("E10").Value = "  +1.21%  "
$ws.Range("E11").Value = "  +1.63%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07856"
$ws.Range("E12").Value = "  -0.37%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.819.85"
$ws.Range("E13").Value = "  +0.07%  "
$ws.Range("E14").Value = "  +0.95%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.520"
$ws.Range("E15").Value = "  -0.45%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "91.13"
$ws.Range("E16").Value = "  -0.34%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.002"
$ws.Range("E17").Value = "  -0.35%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008929"
$ws.Range("E18").Value = "  +1.04%  "
$ws.Range("E19").Value = "  -0.27%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.69"
$ws.Range("E20").Value = "  -0.73%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "27.263.71"
$ws.Range("E21").Value = "  +0.87%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.087"
$ws.Range("E22").Value = "  -0.28%  "
$ws.Range("E23").Value = "  -0.12%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.071.50"
$ws.Range("E24").Value = "  +1.84%  "
$ws.Range("E25").Value = "  +5.29%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "151.31"
$ws.Range("E26").Value = "  -0.97%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.38"
$ws.Range("E27").Value = "  -0.64%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.050"
$ws.Range("E28").Value = "  +0.23%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "115.85"
$ws.Range("E29").Value = "  +0.12%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.046"
$ws.Range("E30").Value = "  -1.53%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08801"
$ws.Range("E31").Value = "  -0.97%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.141"
$ws.Range("E32").Value = "  +6.29%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7680"
$ws.Range("E33").Value = "  +5.44%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.166"
$ws.Range("E34").Value = "  +3.02%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.514"
$ws.Range("E35").Value = "  +1.66%  "
$ws.Range("E36").Value = "  +10.76%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.103"
$ws.Range("E37").Value = "  +2.43%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01937"
$ws.Range("E38").Value = "  -0.67%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05217"
$ws.Range("E39").Value = "  -0.20%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.938"
$ws.Range("E40").Value = "  -0.40%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.042"
$ws.Range("E41").Value = "  -0.91%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5128"
$ws.Range("E42").Value = "  -0.76%  "
$ws.Range("E43").Value = "  +0.18%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.452"
$ws.Range("E44").Value = "  +3.30%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4798"
$ws.Range("E45").Value = "  -0.87%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.40"
$ws.Range("E46").Value = "  +1.70%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.000"
$ws.Range("E47").Value = "  -0.32%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "102.91"
$ws.Range("E48").Value = "  +0.51%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.643"
$ws.Range("E49").Value = "  +0.83%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06206"
$ws.Range("E50").Value = "  +0.03%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "65.53"
$ws.Range("E51").Value = "  +1.13%  "
